$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 218.94737
$ws.Range("I33").Value = 169.375
$ws.Range("J33").Value = 483.33334
$ws.Range("K33").Value = 169.375
$ws.Range("L33").Value = 483.33334
$ws.Range("M33").Value = 59.625
$ws.Range("N33").Value = -941.33334
$ws.Range("H64").Value = 3964.4443
$ws.Range("I64").Value = 3964.4443
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3964.4443
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3716.4443
$ws.Range("H67").Value = 3964.4443
$ws.Range("I67").Value = 3964.4443
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3964.4443
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -3106.4443
$ws.Range("H76").Value = 3025.8064
$ws.Range("I76").Value = 3040
$ws.Range("J76").Value = 2966.6667
$ws.Range("K76").Value = 3040
$ws.Range("L76").Value = 2966.6667
$ws.Range("M76").Value = -2725
$ws.Range("N76").Value = -3596.6667
$ws.Range("H79").Value = 3025.8064
$ws.Range("I79").Value = 3040
$ws.Range("J79").Value = 2966.6667
$ws.Range("K79").Value = 3040
$ws.Range("L79").Value = 2966.6667
$ws.Range("M79").Value = -1948
$ws.Range("N79").Value = -5150.6667
$ws.Range("H129").Value = 1069.7906
$ws.Range("I129").Value = 449
$ws.Range("J129").Value = 1170.4595
$ws.Range("K129").Value = 1347
$ws.Range("L129").Value = 3511.3785
$ws.Range("M129").Value = 3653
$ws.Range("N129").Value = -13511.3785
$ws.Range("H138").Value = 4742.9634
$ws.Range("I138").Value = 3616.261
$ws.Range("J138").Value = 5182.1865
$ws.Range("K138").Value = 10848.783
$ws.Range("L138").Value = 15546.5595
$ws.Range("M138").Value = -5708.782999999999
$ws.Range("N138").Value = -25826.5595

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1051.7693
$ws.Range("I45").Value = 964.6667
$ws.Range("J45").Value = 1126.4286
$ws.Range("K45").Value = 964.6667
$ws.Range("L45").Value = 1126.4286
$ws.Range("M45").Value = -587.6667
$ws.Range("N45").Value = -1880.4286
$ws.Range("H63").Value = 3112.8125
$ws.Range("I63").Value = 2138.8462
$ws.Range("J63").Value = 7333.3335
$ws.Range("K63").Value = 2138.8462
$ws.Range("L63").Value = 7333.3335
$ws.Range("M63").Value = -1452.8462
$ws.Range("N63").Value = -8705.333500000001
$ws.Range("H66").Value = 3112.8125
$ws.Range("I66").Value = 2138.8462
$ws.Range("J66").Value = 7333.3335
$ws.Range("K66").Value = 10694.231
$ws.Range("L66").Value = 36666.6675
$ws.Range("M66").Value = -7262.231
$ws.Range("N66").Value = -43530.6675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1747.8422
$ws.Range("I105").Value = 1680.6
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1680.6
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 66.40000000000009
$ws.Range("N105").Value = -5494
$ws.Range("H134").Value = 2092.4
$ws.Range("I134").Value = 2100.4348
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 6301.3044
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -3766.3044
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3305.328
$ws.Range("I58").Value = 1229.55
$ws.Range("J58").Value = 4317.9023
$ws.Range("K58").Value = 1229.55
$ws.Range("L58").Value = 4317.9023
$ws.Range("M58").Value = -1026.55
$ws.Range("N58").Value = -4723.9023
$ws.Range("H59").Value = 14275.3
$ws.Range("I59").Value = 10500
$ws.Range("J59").Value = 15219.125
$ws.Range("K59").Value = 10500
$ws.Range("L59").Value = 15219.125
$ws.Range("M59").Value = -9355
$ws.Range("N59").Value = -17509.125
$ws.Range("H62").Value = 3166.6667
$ws.Range("I62").Value = 3133.3333
$ws.Range("J62").Value = 3233.3333
$ws.Range("K62").Value = 3133.3333
$ws.Range("L62").Value = 3233.3333
$ws.Range("M62").Value = -2509.3333
$ws.Range("N62").Value = -4481.3333
$ws.Range("H65").Value = 3166.6667
$ws.Range("I65").Value = 3133.3333
$ws.Range("J65").Value = 3233.3333
$ws.Range("K65").Value = 15666.6665
$ws.Range("L65").Value = 16166.6665
$ws.Range("M65").Value = -12546.6665
$ws.Range("N65").Value = -22406.6665
$ws.Range("H68").Value = 18384
$ws.Range("I68").Value = 14268
$ws.Range("J68").Value = 18795.6
$ws.Range("K68").Value = 14268
$ws.Range("L68").Value = 18795.6
$ws.Range("M68").Value = -13519
$ws.Range("N68").Value = -20293.6
$ws.Range("H71").Value = 18384
$ws.Range("I71").Value = 14268
$ws.Range("J71").Value = 18795.6
$ws.Range("K71").Value = 42804
$ws.Range("L71").Value = 56386.8
$ws.Range("M71").Value = -39060
$ws.Range("N71").Value = -63874.8
$ws.Range("H74").Value = 14383.625
$ws.Range("I74").Value = 5185
$ws.Range("J74").Value = 17449.834
$ws.Range("K74").Value = 5185
$ws.Range("L74").Value = 17449.834
$ws.Range("M74").Value = -4311
$ws.Range("N74").Value = -19197.834
$ws.Range("H77").Value = 14383.625
$ws.Range("I77").Value = 5185
$ws.Range("J77").Value = 17449.834
$ws.Range("K77").Value = 15555
$ws.Range("L77").Value = 52349.50199999999
$ws.Range("M77").Value = -11187
$ws.Range("N77").Value = -61085.50199999999
$ws.Range("H134").Value = 2358.76
$ws.Range("I134").Value = 2027.3823
$ws.Range("J134").Value = 3062.9375
$ws.Range("K134").Value = 6082.1469
$ws.Range("L134").Value = 9188.8125
$ws.Range("M134").Value = -3547.1469
$ws.Range("N134").Value = -14258.8125
$ws.Range("H136").Value = 3305.328
$ws.Range("I136").Value = 1229.55
$ws.Range("J136").Value = 4317.9023
$ws.Range("K136").Value = 3688.65
$ws.Range("L136").Value = 12953.7069
$ws.Range("M136").Value = -1138.65
$ws.Range("N136").Value = -18053.7069

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 31535194
$ws.Range("I131").Value = 166681680
$ws.Range("J131").Value = 15153803
$ws.Range("K131").Value = 500045040
$ws.Range("L131").Value = 45461409
$ws.Range("M131").Value = -500040000
$ws.Range("N131").Value = -45471489
$ws.Range("H137").Value = 27331.91
$ws.Range("I137").Value = 3712.8572
$ws.Range("J137").Value = 31682.79
$ws.Range("K137").Value = 11138.5716
$ws.Range("L137").Value = 95048.37
$ws.Range("M137").Value = -6038.571599999999
$ws.Range("N137").Value = -105248.37

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4561.3667
$ws.Range("I70").Value = 4398.4287
$ws.Range("J70").Value = 4703.9375
$ws.Range("K70").Value = 4398.4287
$ws.Range("L70").Value = 4703.9375
$ws.Range("M70").Value = -4128.4287
$ws.Range("N70").Value = -5243.9375
$ws.Range("H73").Value = 4561.3667
$ws.Range("I73").Value = 4398.4287
$ws.Range("J73").Value = 4703.9375
$ws.Range("K73").Value = 4398.4287
$ws.Range("L73").Value = 4703.9375
$ws.Range("M73").Value = -3462.4287
$ws.Range("N73").Value = -6575.9375
$ws.Range("H80").Value = 87582.38
$ws.Range("I80").Value = 3170.8333
$ws.Range("J80").Value = 159935.14
$ws.Range("K80").Value = 3170.8333
$ws.Range("L80").Value = 159935.14
$ws.Range("M80").Value = -2172.8333
$ws.Range("N80").Value = -161931.14
$ws.Range("H83").Value = 87582.38
$ws.Range("I83").Value = 3170.8333
$ws.Range("J83").Value = 159935.14
$ws.Range("K83").Value = 15854.1665
$ws.Range("L83").Value = 799675.7000000001
$ws.Range("M83").Value = -10862.1665
$ws.Range("N83").Value = -809659.7000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3183.3333
$ws.Range("I136").Value = 1664.1428
$ws.Range("J136").Value = 8500.5
$ws.Range("K136").Value = 4992.428400000001
$ws.Range("L136").Value = 25501.5
$ws.Range("M136").Value = -2442.428400000001
$ws.Range("N136").Value = -30601.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3104.8
$ws.Range("I132").Value = 3971.7222
$ws.Range("J132").Value = 2186.8823
$ws.Range("K132").Value = 11915.1666
$ws.Range("L132").Value = 6560.646900000001
$ws.Range("M132").Value = -9385.1666
$ws.Range("N132").Value = -11620.6469
$ws.Range("H136").Value = 21473.285
$ws.Range("I136").Value = 31226.777
$ws.Range("J136").Value = 3917
$ws.Range("K136").Value = 93680.33099999999
$ws.Range("L136").Value = 11751
$ws.Range("M136").Value = -91130.33099999999
$ws.Range("N136").Value = -16851
